$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 173, shifting rows 173:215 down to 174:216.
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row 173 with the new weekly record.
$ws.Cells.Item(173, 1).Value = 3
$ws.Cells.Item(173, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(173, 3).Value = "Coquimbo"
$ws.Cells.Item(173, 4).Value = 44809
$ws.Cells.Item(173, 5).Value = 5
$ws.Cells.Item(173, 6).Value = "Fruta"
$ws.Cells.Item(173, 7).Value = 100101
$ws.Cells.Item(173, 8).Value = "Berries"
$ws.Cells.Item(173, 9).Value = 100101001
$ws.Cells.Item(173, 10).Value = "Arándano (blue)"
$ws.Cells.Item(173, 11).Value = "Sin especificar"
$ws.Cells.Item(173, 12).Value = "Primera"
$ws.Cells.Item(173, 13).Value = 50
$ws.Cells.Item(173, 14).Value = 13000
$ws.Cells.Item(173, 15).Value = 13000
$ws.Cells.Item(173, 16).Value = 13000
$ws.Cells.Item(173, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(173, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(173, 19).Value = 8667
$ws.Cells.Item(173, 20).Value = 1.5
